$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.820.64"
$ws.Range("E2").Value = "  +3.63%  "
$ws.Range("D3").Value = "2.336.58"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'522.21"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").Value = "'135.16"
$ws.Range("E6").Value = "  +4.14%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.538"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").Value = "2.364.78"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("E10").Value = "  +6.78%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("E12").Value = "  +5.40%  "
$ws.Range("D13").Value = "'0.343"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'23.89"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "2.754.16"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "56.888.15"
$ws.Range("E16").Value = "  +3.74%  "
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "2.364.08"
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("D19").Value = "'10.52"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "'324.81"
$ws.Range("E21").Value = "  +5.00%  "
$ws.Range("D22").Value = "'6.57"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "'60.93"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.162"
$ws.Range("E25").Value = "  +6.88%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.990"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'7.91"
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("E28").Value = "  +9.96%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0750"
$ws.Range("E29").Value = "  +5.63%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'170.71"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "'1.71"
$ws.Range("E31").Value = "  +4.02%  "
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").Value = "'18.37"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'0.991"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("D37").Value = "'0.928"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").Value = "'4.05"
$ws.Range("E38").Value = "  +4.79%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +8.80%  "
$ws.Range("D40").Value = "'37.89"
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("D43").Value = "'137.37"
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").Value = "'280.19"
$ws.Range("E44").Value = "  +9.95%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").Value = "'0.0507"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  +10.67%  "

# Reset style on cells that received numeric-looking text so they stay plain (no quotePrefix style)
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
